$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.24598993899883
$ws.Range("D2").Value = 10.79617379259714
$ws.Range("E2").Value = 16.68946058065848
$ws.Range("F2").Value = 31.92513696040787
$ws.Range("G2").Value = 31.48042673494712
$ws.Range("H2").Value = 15.0561446755569
$ws.Range("J2").Value = 11.74327536271102
$ws.Range("K2").Value = 8.516319565688967
$ws.Range("L2").Value = 8.443170875714459
$ws.Range("N2").Value = 20.15448419976025
$ws.Range("O2").Value = 23.26516580848718
$ws.Range("B3").Value = 16.14040274585248
$ws.Range("D3").Value = 10.80494102798977
$ws.Range("E3").Value = 16.73355619278552
$ws.Range("F3").Value = 31.99507709943502
$ws.Range("G3").Value = 31.53788363305038
$ws.Range("H3").Value = 15.09668349588482
$ws.Range("J3").Value = 11.76922460757771
$ws.Range("K3").Value = 8.198550531639043
$ws.Range("L3").Value = 8.399703183412626
$ws.Range("N3").Value = 20.21003056061408
$ws.Range("O3").Value = 23.32873541578058
$ws.Range("B4").Value = 16.07799981834222
$ws.Range("D4").Value = 10.81178869855039
$ws.Range("E4").Value = 16.76258211002052
$ws.Range("F4").Value = 32.04429370988152
$ws.Range("G4").Value = 31.58132946714779
$ws.Range("H4").Value = 15.1236095946423
$ws.Range("J4").Value = 11.78603090666653
$ws.Range("K4").Value = 7.995418919398895
$ws.Range("L4").Value = 8.373807489249984
$ws.Range("N4").Value = 20.245835940964
$ws.Range("O4").Value = 23.37192512257397
$ws.Range("B5").Value = 16.05320277638404
$ws.Range("D5").Value = 10.81494817484469
$ws.Range("E5").Value = 16.7749017324482
$ws.Range("F5").Value = 32.0659258987129
$ws.Range("G5").Value = 31.60108351369511
$ws.Range("H5").Value = 15.13509424420566
$ws.Range("J5").Value = 11.79309980523959
$ws.Range("K5").Value = 7.910700793386235
$ws.Range("L5").Value = 8.36346161467506
$ws.Range("N5").Value = 20.2608554677891
$ws.Range("O5").Value = 23.39056999317776
$ws.Range("B6").Value = 16.04912408464842
$ws.Range("D6").Value = 10.81549510956181
$ws.Range("E6").Value = 16.77697709248862
$ws.Range("F6").Value = 32.0696130491649
$ws.Range("G6").Value = 31.60448730811827
$ws.Range("H6").Value = 15.13703219668202
$ws.Range("J6").Value = 11.79428690686644
$ws.Range("K6").Value = 7.896518639117302
$ws.Range("L6").Value = 8.361756389404555
$ws.Range("N6").Value = 20.26337536544861
$ws.Range("O6").Value = 23.39372903388864
$ws.Range("B7").Value = 16.07766280665938
$ws.Range("D7").Value = 10.8118298133969
$ws.Range("E7").Value = 16.76274626642849
$ws.Range("F7").Value = 32.04457906992796
$ws.Range("G7").Value = 31.58158758474723
$ws.Range("H7").Value = 15.12376240698901
$ws.Range("J7").Value = 11.78612534796488
$ws.Range("K7").Value = 7.99428412966484
$ws.Range("L7").Value = 8.373667114310082
$ws.Range("N7").Value = 20.2460367628775
$ws.Range("O7").Value = 23.37217234488764
$ws.Range("B8").Value = 16.20909308213556
$ws.Range("D8").Value = 10.7988932196502
$ws.Range("E8").Value = 16.70426024590858
$ws.Range("F8").Value = 31.94794957453461
$ws.Range("G8").Value = 31.49854051050398
$ws.Range("H8").Value = 15.06970026040928
$ws.Range("J8").Value = 11.75204174807009
$ws.Range("K8").Value = 8.4084603370591
$ws.Range("L8").Value = 8.428022924191554
$ws.Range("N8").Value = 20.17328450407965
$ws.Range("O8").Value = 23.28622123071678
$ws.Range("B9").Value = 16.48507977313315
$ws.Range("D9").Value = 10.78511147424431
$ws.Range("E9").Value = 16.60501746178508
$ws.Range("F9").Value = 31.80828760608651
$ws.Range("G9").Value = 31.40064359676185
$ws.Range("H9").Value = 14.97981853203023
$ws.Range("J9").Value = 11.69210699363818
$ws.Range("K9").Value = 9.154126403027545
$ws.Range("L9").Value = 8.54058271326875
$ws.Range("N9").Value = 20.04405110024258
$ws.Range("O9").Value = 23.15069361192404
$ws.Range("B10").Value = 16.69756447433578
$ws.Range("D10").Value = 10.78200142428571
$ws.Range("E10").Value = 16.54147668896082
$ws.Range("F10").Value = 31.73611574349779
$ws.Range("G10").Value = 31.36848763569205
$ws.Range("H10").Value = 14.92359935627139
$ws.Range("J10").Value = 11.65224390194163
$ws.Range("K10").Value = 9.658113557316305
$ws.Range("L10").Value = 8.626465048451044
$ws.Range("N10").Value = 19.95721814751908
$ws.Range("O10").Value = 23.07129429732006
$ws.Range("B11").Value = 16.79602606419325
$ws.Range("D11").Value = 10.78209810442488
$ws.Range("E11").Value = 16.51459614807961
$ws.Range("F11").Value = 31.70989953882142
$ws.Range("G11").Value = 31.3625137248282
$ws.Range("H11").Value = 14.90015143023253
$ws.Range("J11").Value = 11.63500691493673
$ws.Range("K11").Value = 9.877319516166615
$ws.Range("L11").Value = 8.666126295602368
$ws.Range("N11").Value = 19.91946180984916
$ws.Range("O11").Value = 23.03956119002365
$ws.Range("B12").Value = 16.83354287827122
$ws.Range("D12").Value = 10.78235091101568
$ws.Range("E12").Value = 16.50470758986929
$ws.Range("F12").Value = 31.70092343631536
$ws.Range("G12").Value = 31.36149605281708
$ws.Range("H12").Value = 14.89157779195206
$ws.Range("J12").Value = 11.62860810135518
$ws.Range("K12").Value = 9.958842981179846
$ws.Range("L12").Value = 8.681221530152042
$ws.Range("N12").Value = 19.90541414750102
$ws.Range("O12").Value = 23.02817586422053
$ws.Range("B13").Value = 16.82545310326815
$ws.Range("D13").Value = 10.78228686842832
$ws.Range("E13").Value = 16.50682435637151
$ws.Range("F13").Value = 31.70281428764839
$ws.Range("G13").Value = 31.36165988477779
$ws.Range("H13").Value = 14.89341069074208
$ws.Range("J13").Value = 11.62998049529041
$ws.Range("K13").Value = 9.941352046120222
$ws.Range("L13").Value = 8.67796724678257
$ws.Range("N13").Value = 19.90842846687629
$ws.Range("O13").Value = 23.03059981144014
$ws.Range("B14").Value = 16.79910811284589
$ws.Range("D14").Value = 10.78211457697147
$ws.Range("E14").Value = 16.51377679141679
$ws.Range("F14").Value = 31.70914200217444
$ws.Range("G14").Value = 31.36240506109633
$ws.Range("H14").Value = 14.89943994868933
$ws.Range("J14").Value = 11.63447790912812
$ws.Range("K14").Value = 9.884056464050557
$ws.Range("L14").Value = 8.667366710097415
$ws.Range("N14").Value = 19.91830109866494
$ws.Range("O14").Value = 23.03861185684591
$ws.Range("B15").Value = 16.78300037969292
$ws.Range("D15").Value = 10.78203716301165
$ws.Range("E15").Value = 16.51807317513315
$ws.Range("F15").Value = 31.71314181004581
$ws.Range("G15").Value = 31.36302356307969
$ws.Range("H15").Value = 14.90317283309644
$ws.Range("J15").Value = 11.63724941885715
$ws.Range("K15").Value = 9.848766820474117
$ws.Range("L15").Value = 8.660883256909056
$ws.Range("N15").Value = 19.92438087820821
$ws.Range("O15").Value = 23.04360169789165
$ws.Range("B16").Value = 16.69116368213275
$ws.Range("D16").Value = 10.78202542403032
$ws.Range("E16").Value = 16.54327408195109
$ws.Range("F16").Value = 31.73796216855894
$ws.Range("G16").Value = 31.36905221332109
$ws.Range("H16").Value = 14.92517450622641
$ws.Range("J16").Value = 11.65338838371729
$ws.Range("K16").Value = 9.643582080098318
$ws.Range("L16").Value = 8.623884267327494
$ws.Range("N16").Value = 19.95972062910877
$ws.Range("O16").Value = 23.073456442931
$ws.Range("B17").Value = 16.63526673880972
$ws.Range("D17").Value = 10.7824045516331
$ws.Range("E17").Value = 16.55925211210626
$ws.Range("F17").Value = 31.7548831663024
$ws.Range("G17").Value = 31.37496743922943
$ws.Range("H17").Value = 14.9392163168649
$ws.Range("J17").Value = 11.66351847135894
$ws.Range("K17").Value = 9.515102048621916
$ws.Range("L17").Value = 8.601332167893425
$ws.Range("N17").Value = 19.98184649037497
$ws.Range("O17").Value = 23.09289519474116
$ws.Range("B18").Value = 16.60328750305311
$ws.Range("D18").Value = 10.78276495545739
$ws.Range("E18").Value = 16.56863284117196
$ws.Range("F18").Value = 31.76523833612412
$ws.Range("G18").Value = 31.37918432805283
$ws.Range("H18").Value = 14.94749295961685
$ws.Range("J18").Value = 11.66942948412589
$ws.Range("K18").Value = 9.440258720441811
$ws.Range("L18").Value = 8.588417146764488
$ws.Range("N18").Value = 19.99473695148129
$ws.Range("O18").Value = 23.10448865535962
$ws.Range("B19").Value = 16.59249012025276
$ws.Range("D19").Value = 10.78291146262361
$ws.Range("E19").Value = 16.57184174692204
$ws.Range("F19").Value = 31.7688513508798
$ws.Range("G19").Value = 31.3807519863182
$ws.Range("H19").Value = 14.95032967311414
$ws.Range("J19").Value = 11.67144537247734
$ws.Range("K19").Value = 9.414757012447803
$ws.Range("L19").Value = 8.584054289522456
$ws.Range("N19").Value = 19.99912968107925
$ws.Range("O19").Value = 23.10848488143591
$ws.Range("B20").Value = 16.64119954253276
$ws.Range("D20").Value = 10.7823494693025
$ws.Range("E20").Value = 16.55753150265221
$ws.Range("F20").Value = 31.75301745210268
$ws.Range("G20").Value = 31.37425344102975
$ws.Range("H20").Value = 14.93770082728609
$ws.Range("J20").Value = 11.66243136930984
$ws.Range("K20").Value = 9.528877068601226
$ws.Range("L20").Value = 8.603727109368982
$ws.Range("N20").Value = 19.97947416260424
$ws.Range("O20").Value = 23.09078317797598
$ws.Range("B21").Value = 16.80684021131579
$ws.Range("D21").Value = 10.78215932523468
$ws.Range("E21").Value = 16.51172681285194
$ws.Range("F21").Value = 31.70725757857203
$ws.Range("G21").Value = 31.36215241272953
$ws.Range("H21").Value = 14.89766071768117
$ws.Range("J21").Value = 11.63315342742312
$ws.Range("K21").Value = 9.900926148406194
$ws.Range("L21").Value = 8.670478343036462
$ws.Range("N21").Value = 19.91539449462506
$ws.Range("O21").Value = 23.03624138690512
$ws.Range("B22").Value = 16.91643268456096
$ws.Range("D22").Value = 10.78329461751577
$ws.Range("E22").Value = 16.48348389386919
$ws.Range("F22").Value = 31.68289620952067
$ws.Range("G22").Value = 31.36149708874335
$ws.Range("H22").Value = 14.87327326268844
$ws.Range("J22").Value = 11.61476712216362
$ws.Range("K22").Value = 10.13541156237851
$ws.Range("L22").Value = 8.714545151915667
$ws.Range("N22").Value = 19.87497055034319
$ws.Range("O22").Value = 23.00427515937323
$ws.Range("B23").Value = 16.85782794962498
$ws.Range("D23").Value = 10.78257384221455
$ws.Range("E23").Value = 16.4984029496153
$ws.Range("F23").Value = 31.69539096516931
$ws.Range("G23").Value = 31.36118336971008
$ws.Range("H23").Value = 14.8861264057976
$ws.Range("J23").Value = 11.62451191943864
$ws.Range("K23").Value = 10.01106688482245
$ws.Range("L23").Value = 8.690988473720314
$ws.Range("N23").Value = 19.896412689425
$ws.Range("O23").Value = 23.02099924362575
$ws.Range("B24").Value = 16.63851683189974
$ws.Range("D24").Value = 10.78237392827075
$ws.Range("E24").Value = 16.55830878353611
$ws.Range("F24").Value = 31.75385898826803
$ws.Range("G24").Value = 31.37457369740941
$ws.Range("H24").Value = 14.93838534522753
$ws.Range("J24").Value = 11.66292257683125
$ws.Range("K24").Value = 9.522652426764781
$ws.Range("L24").Value = 8.60264419841228
$ws.Range("N24").Value = 19.98054616234606
$ws.Range("O24").Value = 23.09173671919044
$ws.Range("B25").Value = 16.40860683409745
$ws.Range("D25").Value = 10.78760401621119
$ws.Range("E25").Value = 16.63021607302392
$ws.Range("F25").Value = 31.84072794342783
$ws.Range("G25").Value = 31.42015242733595
$ws.Range("H25").Value = 15.00240859453559
$ws.Range("J25").Value = 11.70758578898643
$ws.Range("K25").Value = 8.959896530435199
$ws.Range("L25").Value = 8.509539419371032
$ws.Range("N25").Value = 20.07758176825084
$ws.Range("O25").Value = 23.18381753973987
